$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 7-14: Date (serial), Product, Quantity, Invoice Number, Staff
$data = @(
    @(45280, "Yakuza Teriyaki",     7, 2312200004, "Cherry M. Gillego"),
    @(45280, "Chicano Chili",       5, 2312200004, "Cherry M. Gillego"),
    @(45280, "Gangbanger Tuna",     3, 2312200004, "Cherry M. Gillego"),
    @(45284, "Waddup Che&Bac",      5, 2312240005, "Fredrick James Paolo R. Gillego"),
    @(45284, "Cheese Burger",       3, 2312240005, "Fredrick James Paolo R. Gillego"),
    @(45284, "Hardcore Overload ",  6, 2312240005, "Fredrick James Paolo R. Gillego"),
    @(45284, "Waddup Che&Bac",      5, 2312240006, "Josefe Johnatan M. Gillego"),
    @(45284, "Chicano Chili",       5, 2312240007, "Kristina Franchesca M. Gillego")
)

$startRow = 7
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $entry = $data[$i]

    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
}
